$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers for Stanton number calculations ---
$ws.Range("X1").Value = "Theoretical Stanton Number"
$ws.Range("Y1").Value = "Correlated Stanton Number"

# Match the header formatting (bold, bordered, centered) used by the rest of row 1
$ws.Range("W1").Copy()
$ws.Range("X1:Y1").PasteSpecial(-4122)

# Give the new columns sensible widths, matching the style of their neighbours
$ws.Range("X1").ColumnWidth = $ws.Range("W1").ColumnWidth
$ws.Range("Y1").ColumnWidth = $ws.Range("R1").ColumnWidth

# --- Theoretical Stanton Number values ---
$ws.Range("X2").Value = 0.003467963058845684
$ws.Range("X3").Value = 0.004218452465769974
$ws.Range("X4").Value = 0.004373794654161028
$ws.Range("X5").Value = 0.003464024518265738

# --- Correlated Stanton Number values ---
$ws.Range("Y2").Value = 0.002499325357727647
$ws.Range("Y3").Value = 0.002501440173102507
$ws.Range("Y4").Value = 0.002707067273787979
$ws.Range("Y5").Value = 0.002703523672934813
